$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing hours values
$ws.Range("B7").Value = 2
$ws.Range("B8").Value = 2

# Apply same date number format (numFmtId 16, "d-mmm") as the rest of column A
# before assigning values, so Excel doesn't mint a new auto-date style.
$ws.Range("A13:A16").NumberFormat = $ws.Range("A12").NumberFormat

# New rows 13-16
$ws.Range("A13").Value = (Get-Date -Year 2018 -Month 2 -Day 22 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = "R work; data camp class"

$ws.Range("A14").Value = (Get-Date -Year 2018 -Month 2 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "R work"

$ws.Range("A15").Value = (Get-Date -Year 2018 -Month 3 -Day 6 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B15").Value = 3
$ws.Range("C15").Value = "weekly meeting; creating ggplots"

$ws.Range("A16").Value = (Get-Date -Year 2018 -Month 3 -Day 11 -Hour 0 -Minute 0 -Second 0)
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = "ggplots; summarizng data"

# Update selection to match target
$ws.Range("J17").Select()
